$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 from "Y" to "N" (adds a new shared string entry "N")
$ws.Range("C2").Value = "N"

# Update the active selection to C2
$ws.Range("C2").Select()
